# "Generate Report for Archive"
# The localization-status report was regenerated: the outstanding status
# text "Ready for handoff" became "In Translation", and the Status columns
# (which had been autosized to fit the old, longer text) were re-autosized
# to fit the new, shorter text.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# Update the status text in place everywhere it appears.
$wsOverview.Cells.Replace("Ready for handoff", "In Translation")
$wsZhCn.Cells.Replace("Ready for handoff", "In Translation")
$wsDeDe.Cells.Replace("Ready for handoff", "In Translation")

# Shrink the status columns to match the narrower re-fitted width that
# results from the shorter replacement text.
$wsOverview.Columns.Item(5).ColumnWidth = 12.5   # "zh-cn" status column
$wsOverview.Columns.Item(6).ColumnWidth = 12.5   # "de-de" status column
$wsZhCn.Columns.Item(3).ColumnWidth = 12.5        # "Status" column
$wsDeDe.Columns.Item(3).ColumnWidth = 12.5        # "Status" column
